# "cambios de las fracciones" - update the reporting-period dates (Q3 2022 -> Q4 2022)
# on the main "Reporte de Formatos" sheet, adjust the header row height, and
# update the active view/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# --- Row 7 (table header row) grew taller, e.g. due to wrapped header text ---
$ws.Rows.Item(7).RowHeight = 39

# --- Update the reporting period / validation / update dates for data rows 8:20 ---
# Column B: "Fecha de inicio del periodo que se informa" (period start)  01/07/2022 -> 01/10/2022
$ws.Range("B8:B20").Value = 44835
# Column C: "Fecha de término del periodo que se informa" (period end)   30/09/2022 -> 31/12/2022
$ws.Range("C8:C20").Value = 44926
# Column AD: "Fecha de validación" (validation date)                    10/10/2022 -> 10/01/2023
$ws.Range("AD8:AD20").Value = 44936
# Column AE: "Fecha de actualización" (update date)                     10/10/2022 -> 10/01/2023
$ws.Range("AE8:AE20").Value = 44936

# --- Refresh the sheet's active view / selection state ---
$ws.Activate()
$ws.Range("A2:C2").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
